$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.650.82"
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("D3").Value = "1.596.84"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("E4").Value = "  +0.46%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.97"
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("E6").Value = "  -0.39%  "
$ws.Range("E7").Value = "  +0.47%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.0616"
$ws.Range("E8").Value = "  -0.48%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.246"
$ws.Range("E9").Value = "  -0.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.62"
$ws.Range("E10").Value = "  +0.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0847"
$ws.Range("E11").Value = "  +0.72%  "
$ws.Range("D12").Value = "1.822.96"
$ws.Range("E12").Value = "  -0.08%  "
$ws.Range("D13").Value = "1.602.81"
$ws.Range("E13").Value = "  +0.97%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.04"
$ws.Range("E14").Value = "  -0.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.524"
$ws.Range("E15").Value = "  +0.13%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.78"
$ws.Range("E16").Value = "  -0.89%  "
$ws.Range("D17").Value = "26.637.94"
$ws.Range("E17").Value = "  -0.14%  "
$ws.Range("D18").Value = "0.0₃0739"
$ws.Range("E18").Value = "  -2.99%  "
$ws.Range("E19").Value = "  +0.42%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "208.72"
$ws.Range("E20").Value = "  -0.70%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.13"
$ws.Range("E21").Value = "  -0.79%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.28"
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.25"
$ws.Range("E23").Value = "  -2.45%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.95"
$ws.Range("E24").Value = "  +0.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.14"
$ws.Range("E25").Value = "  +0.66%  "
$ws.Range("E26").Value = "  +0.48%  "
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.26"
$ws.Range("E29").Value = "  -0.41%  "
$ws.Range("E30").Value = "  -2.42%  "
$ws.Range("E31").Value = "  -0.31%  "
$ws.Range("E32").Value = "  -0.65%  "
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("E34").Value = "  +18.99%  "
$ws.Range("D35").Value = "1.278.86"
$ws.Range("E35").Value = "  -0.92%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.49"
$ws.Range("E36").Value = "  +1.10%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.598"
$ws.Range("E37").Value = "  -3.58%  "
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.48"
$ws.Range("E38").Value = "  -1.23%  "
$ws.Range("E39").Value = "  -2.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.822"
$ws.Range("E40").Value = "  -0.55%  "
$ws.Range("E41").Value = "  +1.89%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.42"
$ws.Range("E42").Value = "  -0.33%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.777"
$ws.Range("E43").Value = "  -1.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "62.46"
$ws.Range("E44").Value = "  -1.25%  "
$ws.Range("D45").Value = "1.733.22"
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("E46").Value = "  -1.07%  "
$ws.Range("E47").Value = "  -0.64%  "
$ws.Range("E48").Value = "  +2.12%  "
$ws.Range("E49").Value = "  +0.50%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.49"
$ws.Range("E50").Value = "  +1.78%  "
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.01"
$ws.Range("E51").Value = "  +0.43%  "
